$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("getDataGraphQL")

# Update header cell C1 from "graphQLSentence" to "query"
$ws.Range("C1").Value = "query"

# Update the selected cell/range on the sheet
$ws.Range("C10").Select()
